# Quarterly indexing esoteric bug-fix operation
#
# Column A holds dates that originally mark the 1st day of each calendar
# quarter (e.g. 1988-07-01). They need to be shifted to the 15th of the
# following month (e.g. 1988-08-15) to correct the quarterly indexing bug.
# Column B (growth rates) and everything else on the sheet stays untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 150; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $serial = $cell.Value2
    $d = [DateTime]::FromOADate($serial)
    $shifted = $d.AddMonths(1).AddDays(14)
    $cell.Value = $shifted.ToOADate()
}
